$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add three new header cells (L1:N1), copying the style of the existing
#     header cell K1 so they get the same bold/centered/bordered formatting ---
$ws.Range("K1").Copy($ws.Range("L1:N1"))
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

# --- Update E and F columns (particip / taxa_sucesso) which now store the
#     percentage value already multiplied by 100 instead of a 0-1 fraction ---
$ws.Range("E2").Value = 89.51310861423221
$ws.Range("F2").Value = 61.33891213389121

$ws.Range("E3").Value = 10.48689138576779
$ws.Range("F3").Value = 69.28571428571428

$ws.Range("E4").Value = 92.8474114441417
$ws.Range("F4").Value = 94.13059427732942

$ws.Range("E5").Value = 7.152588555858311
$ws.Range("F5").Value = 95.23809523809523

$ws.Range("E6").Value = 88.74269005847954
$ws.Range("F6").Value = 22.07578253706755

$ws.Range("E7").Value = 11.25730994152047
$ws.Range("F7").Value = 23.37662337662337

# --- Fill in the three new data columns (L: apoio_medio, M: contribuicoes,
#     N: media_contribuicoes) for every data row ---
$ws.Range("L2").Value = 88.07870613099213
$ws.Range("M2").Value = 225303
$ws.Range("N2").Value = 307.3710777626194

$ws.Range("L3").Value = 110.2975973828001
$ws.Range("M3").Value = 38250
$ws.Range("N3").Value = 394.3298969072165

$ws.Range("L4").Value = 88.66083985762999
$ws.Range("M4").Value = 187667
$ws.Range("N4").Value = 146.2720187061574

$ws.Range("L5").Value = 107.8551914385913
$ws.Range("M5").Value = 15979
$ws.Range("N5").Value = 159.79

$ws.Range("L6").Value = 17.98549503340952
$ws.Range("M6").Value = 1940
$ws.Range("N6").Value = 14.47761194029851

$ws.Range("L7").Value = 30.9518559327251
$ws.Range("M7").Value = 268
$ws.Range("N7").Value = 14.88888888888889
